$d = $word.ActiveDocument

# 1. Remove the stray _GoBack bookmark from its old location
#    (it will be re-added later, right before the new "Output Parameter" paragraph)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete() | Out-Null
}

# 2. Append a note to the end of the "String.IsNullOrEmpty(str);" paragraph
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endOfText = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$endOfText.Collapse(0)
$endOfText.InsertAfter("   (checks for null like option 3)") | Out-Null

# 3. Append the new paragraphs (ReadLine/String.Format/String.Interpolation/TryParse
#    notes plus the new "Parameter Kinds" section) at the end of the document body.
$contentEnd = $d.Content
$contentEnd.Collapse(0)
$newParagraphsXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>ReadLine</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>) returns a string</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>String.Format</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>String.Interpolation</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p/><w:p><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>TryParse</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">string) </w:t></w:r><w:r><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:t xml:space="preserve"> evaluates and returns false if not able to parse</w:t></w:r></w:p><w:p/><w:p/><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Parameter Kinds</w:t></w:r></w:p><w:p><w:r><w:t>Input Parameter</w:t></w:r><w:r><w:t xml:space="preserve"> [</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>param</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">] – also referred to as ‘passed by value’. Doesn’t change original variable. </w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Input/Output</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Parameter [ref </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>param</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>] – passed by reference.  Changes original variable (rarely used)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Output Parameter </w:t></w:r></w:p>'
$contentEnd.InsertXML($newParagraphsXml) | Out-Null
